$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.487.51'
$ws.Range("E2").Value = '  +0.50%  '

$ws.Range("D3").Value = '1.808.69'
$ws.Range("E3").Value = '  +0.64%  '

$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").Value = '224.77'
$ws.Range("E5").Value = '  -0.92%  '

$ws.Range("D6").Value = '0.605'
$ws.Range("E6").Value = '  +5.73%  '

$ws.Range("E7").Value = '  -0.10%  '

$ws.Range("D8").Value = '38.28'
$ws.Range("E8").Value = '  +6.21%  '

$ws.Range("E9").Value = '  -3.87%  '

$ws.Range("E10").Value = '  -2.83%  '

$ws.Range("D11").Value = '0.0980'
$ws.Range("E11").Value = '  +2.10%  '

$ws.Range("D12").Value = '2.067.78'
$ws.Range("E12").Value = '  +0.54%  '

$ws.Range("D13").Value = '11.16'
$ws.Range("E13").Value = '  -2.64%  '

$ws.Range("D14").Value = '1.806.89'
$ws.Range("E14").Value = '  +0.46%  '

$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").Value = '0.630'
$ws.Range("E15").Value = '  -1.44%  '

$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '34.451.78'
$ws.Range("E16").Value = '  +0.48%  '

$ws.Range("D17").Value = '4.38'
$ws.Range("E17").Value = '  -2.15%  '

$ws.Range("D18").Value = '68.25'
$ws.Range("E18").Value = '  -0.62%  '

$ws.Range("D19").Value = '241.14'
$ws.Range("E19").Value = '  -1.02%  '

$ws.Range("D20").Value = '0.0₃0768'
$ws.Range("E20").Value = '  -2.49%  '

$ws.Range("D21").Value = '11.14'
$ws.Range("E21").Value = '  -3.46%  '

$ws.Range("E22").Value = '  -0.09%  '

$ws.Range("D23").Value = '4.10'
$ws.Range("E23").Value = '  -1.11%  '

$ws.Range("E24").Value = '  +1.36%  '

$ws.Range("D25").Value = '170.70'
$ws.Range("E25").Value = '  -0.77%  '

$ws.Range("D26").Value = '7.66'
$ws.Range("E26").Value = '  -3.08%  '

$ws.Range("E27").Value = '  +4.32%  '

$ws.Range("D28").Value = '0.122'
$ws.Range("E28").Value = '  +4.09%  '

$ws.Range("E29").Value = '  -0.06%  '

$ws.Range("E30").Value = '  -0.93%  '

$ws.Range("D31").Value = '3.77'
$ws.Range("E31").Value = '  -1.02%  '

$ws.Range("E32").Value = '  -2.15%  '

$ws.Range("D33").Value = '3.84'
$ws.Range("E33").Value = '  -3.57%  '

$ws.Range("E34").Value = '  +1.84%  '

$ws.Range("B35").Value = 'Maker'
$ws.Range("C35").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D35").Value = '1.311.08'
$ws.Range("E35").Value = '  -5.63%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '0.640'
$ws.Range("E36").Value = '  -3.86%  '

$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '1.05'
$ws.Range("E37").Value = '  -0.55%  '

$ws.Range("E38").Value = '  -1.19%  '

$ws.Range("E39").Value = '  -5.16%  '

$ws.Range("B40").Value = 'WEMIXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D40").Value = '1.22'
$ws.Range("E40").Value = '  +3.45%  '

$ws.Range("B41").Value = 'HuobiToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D41").Value = '2.44'
$ws.Range("E41").Value = '  +1.48%  '

$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = '82.28'
$ws.Range("E42").Value = '  +0.22%  '

$ws.Range("E43").Value = '  +0.07%  '

$ws.Range("D44").Value = '0.950'
$ws.Range("E44").Value = '  -0.60%  '

$ws.Range("D45").Value = '13.98'
$ws.Range("E45").Value = '  +5.45%  '

$ws.Range("E46").Value = '  +2.50%  '

$ws.Range("D47").Value = '1.968.61'
$ws.Range("E47").Value = '  +0.57%  '

$ws.Range("E48").Value = '  -3.36%  '

$ws.Range("E49").Value = '  -0.06%  '

$ws.Range("D50").Value = '102.75'
$ws.Range("E50").Value = '  -0.94%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.0609'
$ws.Range("E51").Value = '  +0.29%  '
